# Revert "Revert "Update azure-pipelines.yml"" — restores the D010/D012
# defect rows to their earlier (pre-"Revert") content on the Defects sheet:
#   - D010 (row 11): Status goes back from "Rejected" to "Open", and the
#     Coments cell explaining the rejection is cleared out entirely.
#   - D012 (row 13): the Description and Status cells are cleared out,
#     leaving only the Defect ID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D010 row: Status -> "Open", drop the "Coments" cell completely.
$ws.Range("C11").Value = "Open"
$ws.Range("D11").Clear()

# D012 row: drop Description + Status cells completely.
$ws.Range("B13").Clear()
$ws.Range("C13").Clear()

# Leave the selection the way the author left it: the whole of row 9.
$ws.Rows(9).Select()
